$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the e8dffba3-... file. Mark it "Ready for handoff". ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-12 12:56:43"

# --- zh-cn sheet: row 3 is the e8dffba3-... file. ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-12 12:56:36"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/fc31ca5669e9b2a14343778fb18c820fd2368487/e2e/e8dffba3-cc15-45e5-9cdc-d6048bd3d3f5.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/52ae2005130ffd42468d36f7aaa5ac2ef82e8d0e/e2e/e8dffba3-cc15-45e5-9cdc-d6048bd3d3f5.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 is the e8dffba3-... file. ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-12 12:56:43"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/fc31ca5669e9b2a14343778fb18c820fd2368487/e2e/e8dffba3-cc15-45e5-9cdc-d6048bd3d3f5.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/52ae2005130ffd42468d36f7aaa5ac2ef82e8d0e/e2e/e8dffba3-cc15-45e5-9cdc-d6048bd3d3f5.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
